$wb = $excel.ActiveWorkbook

# Activate the "miasta" worksheet (third sheet)
$ws = $wb.Worksheets.Item("miasta")
$ws.Activate()

# Fix the erroneous city name and its wage value in row 6
$ws.Range("A6").Value = "Rzeszów"
$ws.Range("B6").Value = 7187.74

# Leave selection on A7 to match the final saved state
$ws.Range("A7").Select()
